$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "버스"
$ws.Cells.Item(2, 3).Value = 5925
$ws.Cells.Item(3, 2).Value = "노선"
$ws.Cells.Item(3, 3).Value = 1946
$ws.Cells.Item(4, 2).Value = "출근"
$ws.Cells.Item(4, 3).Value = 1911
$ws.Cells.Item(5, 2).Value = "시간"
$ws.Cells.Item(5, 3).Value = 1640
$ws.Cells.Item(6, 2).Value = "퇴근"
$ws.Cells.Item(6, 3).Value = 1582
$ws.Cells.Item(7, 2).Value = "선"
$ws.Cells.Item(7, 3).Value = 1305
$ws.Cells.Item(8, 2).Value = "없"
$ws.Cells.Item(8, 3).Value = 1206
$ws.Cells.Item(9, 2).Value = "분"
$ws.Cells.Item(9, 3).Value = 1178
$ws.Cells.Item(10, 2).Value = "출발"
$ws.Cells.Item(10, 3).Value = 1080
$ws.Cells.Item(11, 2).Value = "운행"
$ws.Cells.Item(11, 3).Value = 1018
$ws.Cells.Item(12, 2).Value = "셔틀"
$ws.Cells.Item(12, 3).Value = 1005
$ws.Cells.Item(13, 2).Value = "통근"
$ws.Cells.Item(13, 3).Value = 988
$ws.Cells.Item(14, 2).Value = "정류장"
$ws.Cells.Item(14, 3).Value = 984
$ws.Cells.Item(15, 2).Value = "이용"
$ws.Cells.Item(15, 3).Value = 972
$ws.Cells.Item(16, 2).Value = "사람"
$ws.Cells.Item(16, 3).Value = 833
$ws.Cells.Item(17, 2).Value = "많"
$ws.Cells.Item(17, 3).Value = 821
$ws.Cells.Item(18, 2).Value = "부탁"
$ws.Cells.Item(18, 3).Value = 778
$ws.Cells.Item(19, 2).Value = "도착"
$ws.Cells.Item(19, 3).Value = 754
$ws.Cells.Item(20, 2).Value = "시"
$ws.Cells.Item(20, 3).Value = 707
$ws.Cells.Item(21, 2).Value = "역"
$ws.Cells.Item(21, 3).Value = 677
$ws.Cells.Item(22, 2).Value = "앞"
$ws.Cells.Item(22, 3).Value = 662
$ws.Cells.Item(23, 2).Value = "기사님"
$ws.Cells.Item(23, 3).Value = 636
$ws.Cells.Item(24, 2).Value = "차량"
$ws.Cells.Item(24, 3).Value = 620
$ws.Cells.Item(25, 2).Value = "정차"
$ws.Cells.Item(25, 3).Value = 600
$ws.Cells.Item(26, 2).Value = "차"
$ws.Cells.Item(26, 3).Value = 574
$ws.Cells.Item(27, 2).Value = "감사"
$ws.Cells.Item(27, 3).Value = 569
$ws.Cells.Item(28, 2).Value = "요청"
$ws.Cells.Item(28, 3).Value = 561
$ws.Cells.Item(29, 2).Value = "생각"
$ws.Cells.Item(29, 3).Value = 558
$ws.Cells.Item(30, 2).Value = "출퇴근"
$ws.Cells.Item(30, 3).Value = 554
$ws.Cells.Item(31, 2).Value = "회사"
$ws.Cells.Item(31, 3).Value = 542
$ws.Cells.Item(32, 2).Value = "탑승"
$ws.Cells.Item(32, 3).Value = 513
$ws.Cells.Item(33, 2).Value = "변경"
$ws.Cells.Item(33, 3).Value = 501
$ws.Cells.Item(34, 2).Value = "인원"
$ws.Cells.Item(34, 3).Value = 477
$ws.Cells.Item(35, 2).Value = "아침"
$ws.Cells.Item(35, 3).Value = 477
$ws.Cells.Item(36, 2).Value = "경우"
$ws.Cells.Item(36, 3).Value = 473
$ws.Cells.Item(37, 2).Value = "문의"
$ws.Cells.Item(37, 3).Value = 465
$ws.Cells.Item(38, 2).Value = "관련"
$ws.Cells.Item(38, 3).Value = 462
$ws.Cells.Item(39, 2).Value = "검토"
$ws.Cells.Item(39, 3).Value = 454
$ws.Cells.Item(40, 2).Value = "정도"
$ws.Cells.Item(40, 3).Value = 450
$ws.Cells.Item(41, 2).Value = "청주"
$ws.Cells.Item(41, 3).Value = 450
$ws.Cells.Item(42, 2).Value = "동백"
$ws.Cells.Item(42, 3).Value = 450
$ws.Cells.Item(43, 2).Value = "좋"
$ws.Cells.Item(43, 3).Value = 430
$ws.Cells.Item(44, 2).Value = "안녕"
$ws.Cells.Item(44, 3).Value = 429
$ws.Cells.Item(45, 2).Value = "건의"
$ws.Cells.Item(45, 3).Value = 427
$ws.Cells.Item(46, 2).Value = "불편"
$ws.Cells.Item(46, 3).Value = 424
$ws.Cells.Item(47, 2).Value = "이천"
$ws.Cells.Item(47, 3).Value = 421
$ws.Cells.Item(48, 2).Value = "신설"
$ws.Cells.Item(48, 3).Value = 419
$ws.Cells.Item(49, 2).Value = "답변"
$ws.Cells.Item(49, 3).Value = 399
$ws.Cells.Item(50, 2).Value = "수원"
$ws.Cells.Item(50, 3).Value = 394
$ws.Cells.Item(51, 2).Value = "사항"
$ws.Cells.Item(51, 3).Value = 378
$ws.Cells.Item(52, 2).Value = "추가"
$ws.Cells.Item(52, 3).Value = 374
$ws.Cells.Item(53, 2).Value = "조"
$ws.Cells.Item(53, 3).Value = 368
$ws.Cells.Item(54, 2).Value = "위치"
$ws.Cells.Item(54, 3).Value = 367
$ws.Cells.Item(55, 2).Value = "교대"
$ws.Cells.Item(55, 3).Value = 366
$ws.Cells.Item(56, 2).Value = "용인"
$ws.Cells.Item(56, 3).Value = 360
$ws.Cells.Item(57, 2).Value = "글"
$ws.Cells.Item(57, 3).Value = 354
$ws.Cells.Item(58, 2).Value = "오늘"
$ws.Cells.Item(58, 3).Value = 340
$ws.Cells.Item(59, 2).Value = "동"
$ws.Cells.Item(59, 3).Value = 338
$ws.Cells.Item(60, 2).Value = "전"
$ws.Cells.Item(60, 3).Value = 336
$ws.Cells.Item(61, 2).Value = "중"
$ws.Cells.Item(61, 3).Value = 325
$ws.Cells.Item(62, 2).Value = "아파트"
$ws.Cells.Item(62, 3).Value = 321
$ws.Cells.Item(63, 2).Value = "일"
$ws.Cells.Item(63, 3).Value = 308
$ws.Cells.Item(64, 2).Value = "현재"
$ws.Cells.Item(64, 3).Value = 291
$ws.Cells.Item(65, 2).Value = "신"
$ws.Cells.Item(65, 3).Value = 290
$ws.Cells.Item(66, 2).Value = "쪽"
$ws.Cells.Item(66, 3).Value = 289
$ws.Cells.Item(67, 2).Value = "수고"
$ws.Cells.Item(67, 3).Value = 289
$ws.Cells.Item(68, 2).Value = "늦"
$ws.Cells.Item(68, 3).Value = 286
$ws.Cells.Item(69, 2).Value = "하차"
$ws.Cells.Item(69, 3).Value = 285
$ws.Cells.Item(70, 2).Value = "곳"
$ws.Cells.Item(70, 3).Value = 277
$ws.Cells.Item(71, 2).Value = "후"
$ws.Cells.Item(71, 3).Value = 273
$ws.Cells.Item(72, 2).Value = "확인"
$ws.Cells.Item(72, 3).Value = 265
$ws.Cells.Item(73, 2).Value = "직원"
$ws.Cells.Item(73, 3).Value = 265
$ws.Cells.Item(74, 2).Value = "인천"
$ws.Cells.Item(74, 3).Value = 259
$ws.Cells.Item(75, 2).Value = "안"
$ws.Cells.Item(75, 3).Value = 257
$ws.Cells.Item(76, 2).Value = "때문"
$ws.Cells.Item(76, 3).Value = 253
$ws.Cells.Item(77, 2).Value = "운전"
$ws.Cells.Item(77, 3).Value = 250
$ws.Cells.Item(78, 2).Value = "지역"
$ws.Cells.Item(78, 3).Value = 248
$ws.Cells.Item(79, 2).Value = "기사분"
$ws.Cells.Item(79, 3).Value = 238
$ws.Cells.Item(80, 2).Value = "운영"
$ws.Cells.Item(80, 3).Value = 238
$ws.Cells.Item(81, 2).Value = "공지"
$ws.Cells.Item(81, 3).Value = 236
$ws.Cells.Item(82, 2).Value = "춥"
$ws.Cells.Item(82, 3).Value = 234
$ws.Cells.Item(83, 2).Value = "서울"
$ws.Cells.Item(83, 3).Value = 230
$ws.Cells.Item(84, 2).Value = "안녕하"
$ws.Cells.Item(84, 3).Value = 225
$ws.Cells.Item(85, 2).Value = "배차"
$ws.Cells.Item(85, 3).Value = 219
$ws.Cells.Item(86, 2).Value = "문제"
$ws.Cells.Item(86, 3).Value = 217
$ws.Cells.Item(87, 2).Value = "길"
$ws.Cells.Item(87, 3).Value = 217
$ws.Cells.Item(88, 2).Value = "히터"
$ws.Cells.Item(88, 3).Value = 215
$ws.Cells.Item(89, 2).Value = "하세"
$ws.Cells.Item(89, 3).Value = 210
$ws.Cells.Item(90, 2).Value = "대"
$ws.Cells.Item(90, 3).Value = 209
$ws.Cells.Item(91, 2).Value = "개선"
$ws.Cells.Item(91, 3).Value = 208
$ws.Cells.Item(92, 2).Value = "주세"
$ws.Cells.Item(92, 3).Value = 205
$ws.Cells.Item(93, 2).Value = "거리"
$ws.Cells.Item(93, 3).Value = 204
$ws.Cells.Item(94, 2).Value = "조사"
$ws.Cells.Item(94, 3).Value = 203
$ws.Cells.Item(95, 2).Value = "한번"
$ws.Cells.Item(95, 3).Value = 202
$ws.Cells.Item(96, 2).Value = "증차"
$ws.Cells.Item(96, 3).Value = 200
$ws.Cells.Item(97, 2).Value = "데"
$ws.Cells.Item(97, 3).Value = 199
$ws.Cells.Item(98, 2).Value = "잠실"
$ws.Cells.Item(98, 3).Value = 196
$ws.Cells.Item(99, 2).Value = "자리"
$ws.Cells.Item(99, 3).Value = 195
$ws.Cells.Item(100, 2).Value = "담당자"
$ws.Cells.Item(100, 3).Value = 194
$ws.Cells.Item(101, 2).Value = "문"
$ws.Cells.Item(101, 3).Value = 193
